$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-10 06:34:55"

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
